$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- BB1 header cell: copy style from BA1 (date header), then set new date value ---
$srcHeader = $ws.Cells.Item(1,53)
$dstHeader = $ws.Cells.Item(1,54)
$srcHeader.Copy($dstHeader)
$dstHeader.Value = 45986

# --- BB2:BB82 data cells (no special style, same as BA column) ---
$ws.Cells.Item(2,54).Value = -1
$ws.Cells.Item(3,54).Value = 0.5
$ws.Cells.Item(4,54).Value = -0.4
$ws.Cells.Item(5,54).Value = 0.2
$ws.Cells.Item(6,54).Value = -0.2
$ws.Cells.Item(7,54).Value = 0
$ws.Cells.Item(8,54).Value = -1.6
$ws.Cells.Item(9,54).Value = 1.2
$ws.Cells.Item(10,54).Value = -0.5
$ws.Cells.Item(11,54).Value = 1
$ws.Cells.Item(12,54).Value = 0.3
$ws.Cells.Item(13,54).Value = -0.4
$ws.Cells.Item(14,54).Value = 0.2
$ws.Cells.Item(15,54).Value = -0.3
$ws.Cells.Item(16,54).Value = 1.2
$ws.Cells.Item(17,54).Value = -0.7
$ws.Cells.Item(18,54).Value = -0.2
$ws.Cells.Item(19,54).Value = -0.3
$ws.Cells.Item(20,54).Value = -0.6
$ws.Cells.Item(21,54).Value = -0.4
$ws.Cells.Item(22,54).Value = 0.5
$ws.Cells.Item(23,54).Value = 0
$ws.Cells.Item(24,54).Value = -0.1
$ws.Cells.Item(25,54).Value = 0.3
$ws.Cells.Item(26,54).Value = 0.1
$ws.Cells.Item(27,54).Value = 0
$ws.Cells.Item(28,54).Value = 0.4
$ws.Cells.Item(29,54).Value = -0.7
$ws.Cells.Item(30,54).Value = 0
$ws.Cells.Item(31,54).Value = -0.1
$ws.Cells.Item(32,54).Value = -0.5
$ws.Cells.Item(33,54).Value = 0.2
$ws.Cells.Item(34,54).Value = 0.1
$ws.Cells.Item(35,54).Value = -0.4
$ws.Cells.Item(36,54).Value = 0
$ws.Cells.Item(37,54).Value = 0.3
$ws.Cells.Item(38,54).Value = 0.4
$ws.Cells.Item(39,54).Value = -0.5
$ws.Cells.Item(40,54).Value = 0.2
$ws.Cells.Item(41,54).Value = 0.1
$ws.Cells.Item(42,54).Value = -0.1
$ws.Cells.Item(43,54).Value = -0.1
$ws.Cells.Item(44,54).Value = 0.3
$ws.Cells.Item(45,54).Value = 0.4
$ws.Cells.Item(46,54).Value = -0.1
$ws.Cells.Item(47,54).Value = -0.8
$ws.Cells.Item(48,54).Value = 0.2
$ws.Cells.Item(49,54).Value = -0.9
$ws.Cells.Item(50,54).Value = 0.3
$ws.Cells.Item(51,54).Value = 0.1
$ws.Cells.Item(52,54).Value = -0.1
$ws.Cells.Item(53,54).Value = -1.9
$ws.Cells.Item(54,54).Value = 0.5
$ws.Cells.Item(55,54).Value = 2.2
$ws.Cells.Item(56,54).Value = -0.9
$ws.Cells.Item(57,54).Value = 0.1
$ws.Cells.Item(58,54).Value = 0.2
$ws.Cells.Item(59,54).Value = -0.3
$ws.Cells.Item(60,54).Value = 0.7
$ws.Cells.Item(61,54).Value = 0.2
$ws.Cells.Item(62,54).Value = 0.4
$ws.Cells.Item(63,54).Value = -0.1
$ws.Cells.Item(64,54).Value = 0.2
$ws.Cells.Item(65,54).Value = -0.6
$ws.Cells.Item(66,54).Value = -0.7
$ws.Cells.Item(67,54).Value = 0
$ws.Cells.Item(68,54).Value = 0
$ws.Cells.Item(69,54).Value = 0.9
$ws.Cells.Item(70,54).Value = 0.5
$ws.Cells.Item(71,54).Value = -0.5
$ws.Cells.Item(72,54).Value = 0.4
$ws.Cells.Item(73,54).Value = 0.2
$ws.Cells.Item(74,54).Value = -0.166915309931608
$ws.Cells.Item(75,54).Value = -0.0273035672071738
$ws.Cells.Item(76,54).Value = 0.06878353535742468
$ws.Cells.Item(77,54).Value = 0.00284077231157528
$ws.Cells.Item(78,54).Value = -0.01653670234748227
$ws.Cells.Item(79,54).Value = 0.009488806423010147
$ws.Cells.Item(80,54).Value = 0.01089042766875424
$ws.Cells.Item(81,54).Value = 0.00184892501921626
$ws.Cells.Item(82,54).Value = 0.003157254381447184

# --- New row 83: A83 date cell (copy style from A82), BB83 data cell ---
$srcA = $ws.Cells.Item(82,1)
$dstA = $ws.Cells.Item(83,1)
$srcA.Copy($dstA)
$dstA.Value = 46934
$ws.Cells.Item(83,54).Value = 0.005944842580863992
